$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44187
$ws.Cells.Item(2, 11).Value = 'Dina'
$ws.Cells.Item(2, 13).Value = 240
$ws.Cells.Item(2, 14).Value = 22000
$ws.Cells.Item(2, 15).Value = 23000
$ws.Cells.Item(2, 16).Value = 22500
$ws.Cells.Item(2, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(2, 19).Value = 1250
$ws.Cells.Item(2, 20).Value = 18

# Row 3
$ws.Cells.Item(3, 4).Value = 44174
$ws.Cells.Item(3, 13).Value = 240
$ws.Cells.Item(3, 14).Value = 22500
$ws.Cells.Item(3, 15).Value = 23000
$ws.Cells.Item(3, 16).Value = 22750
$ws.Cells.Item(3, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(3, 19).Value = 1264
$ws.Cells.Item(3, 20).Value = 18

# Row 4
$ws.Cells.Item(4, 4).Value = 44165
$ws.Cells.Item(4, 11).Value = 'Castle Brite'
$ws.Cells.Item(4, 12).Value = 'Especial'
$ws.Cells.Item(4, 13).Value = 240
$ws.Cells.Item(4, 14).Value = 20500
$ws.Cells.Item(4, 15).Value = 21000
$ws.Cells.Item(4, 16).Value = 20750
$ws.Cells.Item(4, 17).Value = '$/caja 15 kilos'
$ws.Cells.Item(4, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(4, 19).Value = 1383
$ws.Cells.Item(4, 20).Value = 15

# Row 5
$ws.Cells.Item(5, 4).Value = 44165
$ws.Cells.Item(5, 11).Value = 'Castle Brite'
$ws.Cells.Item(5, 12).Value = 'Primera'
$ws.Cells.Item(5, 14).Value = 17500
$ws.Cells.Item(5, 15).Value = 18000
$ws.Cells.Item(5, 16).Value = 17750
$ws.Cells.Item(5, 17).Value = '$/caja 15 kilos'
$ws.Cells.Item(5, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(5, 19).Value = 1183
$ws.Cells.Item(5, 20).Value = 15

# Row 6
$ws.Cells.Item(6, 4).Value = 44167
$ws.Cells.Item(6, 11).Value = 'Castle Brite'
$ws.Cells.Item(6, 12).Value = 'Especial'
$ws.Cells.Item(6, 13).Value = 400
$ws.Cells.Item(6, 14).Value = 20000
$ws.Cells.Item(6, 15).Value = 21000
$ws.Cells.Item(6, 16).Value = 20500
$ws.Cells.Item(6, 17).Value = '$/caja 15 kilos'
$ws.Cells.Item(6, 19).Value = 1367
$ws.Cells.Item(6, 20).Value = 15

# Row 7
$ws.Cells.Item(7, 4).Value = 44167
$ws.Cells.Item(7, 11).Value = 'Castle Brite'
$ws.Cells.Item(7, 12).Value = 'Primera'
$ws.Cells.Item(7, 13).Value = 360
$ws.Cells.Item(7, 14).Value = 17000
$ws.Cells.Item(7, 15).Value = 18000
$ws.Cells.Item(7, 16).Value = 17500
$ws.Cells.Item(7, 17).Value = '$/caja 15 kilos'
$ws.Cells.Item(7, 19).Value = 1167
$ws.Cells.Item(7, 20).Value = 15

# Row 8
$ws.Cells.Item(8, 4).Value = 44161
$ws.Cells.Item(8, 11).Value = 'Dina'
$ws.Cells.Item(8, 12).Value = 'Primera'
$ws.Cells.Item(8, 13).Value = 240
$ws.Cells.Item(8, 14).Value = 19500
$ws.Cells.Item(8, 15).Value = 20000
$ws.Cells.Item(8, 16).Value = 19750
$ws.Cells.Item(8, 19).Value = 1317

# Row 9
$ws.Cells.Item(9, 4).Value = 44161
$ws.Cells.Item(9, 11).Value = 'Dina'
$ws.Cells.Item(9, 12).Value = 'Segunda'
$ws.Cells.Item(9, 13).Value = 140
$ws.Cells.Item(9, 14).Value = 17500
$ws.Cells.Item(9, 16).Value = 17750
$ws.Cells.Item(9, 19).Value = 1183

# Row 10
$ws.Cells.Item(10, 4).Value = 44162
$ws.Cells.Item(10, 12).Value = 'Especial'
$ws.Cells.Item(10, 13).Value = 300
$ws.Cells.Item(10, 14).Value = 20500
$ws.Cells.Item(10, 15).Value = 21000
$ws.Cells.Item(10, 16).Value = 20750
$ws.Cells.Item(10, 17).Value = '$/caja 15 kilos'
$ws.Cells.Item(10, 19).Value = 1383
$ws.Cells.Item(10, 20).Value = 15

# Row 11
$ws.Cells.Item(11, 4).Value = 44162
$ws.Cells.Item(11, 13).Value = 300
$ws.Cells.Item(11, 14).Value = 17500
$ws.Cells.Item(11, 15).Value = 18000
$ws.Cells.Item(11, 16).Value = 17750
$ws.Cells.Item(11, 17).Value = '$/caja 15 kilos'
$ws.Cells.Item(11, 19).Value = 1183
$ws.Cells.Item(11, 20).Value = 15

# Row 12
$ws.Cells.Item(12, 4).Value = 44189
$ws.Cells.Item(12, 13).Value = 400
$ws.Cells.Item(12, 14).Value = 23500
$ws.Cells.Item(12, 15).Value = 24000
$ws.Cells.Item(12, 16).Value = 23750
$ws.Cells.Item(12, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(12, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(12, 19).Value = 1319
$ws.Cells.Item(12, 20).Value = 18

# Row 13
$ws.Cells.Item(13, 4).Value = 44189
$ws.Cells.Item(13, 13).Value = 200
$ws.Cells.Item(13, 14).Value = 21500
$ws.Cells.Item(13, 15).Value = 22000
$ws.Cells.Item(13, 16).Value = 21750
$ws.Cells.Item(13, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(13, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(13, 19).Value = 1208
$ws.Cells.Item(13, 20).Value = 18

# Row 14
$ws.Cells.Item(14, 4).Value = 44186
$ws.Cells.Item(14, 11).Value = 'Dina'
$ws.Cells.Item(14, 13).Value = 200
$ws.Cells.Item(14, 14).Value = 22500
$ws.Cells.Item(14, 15).Value = 23000
$ws.Cells.Item(14, 16).Value = 22750
$ws.Cells.Item(14, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(14, 19).Value = 1264
$ws.Cells.Item(14, 20).Value = 18

# Row 15
$ws.Cells.Item(15, 4).Value = 44168
$ws.Cells.Item(15, 14).Value = 23500
$ws.Cells.Item(15, 15).Value = 24000
$ws.Cells.Item(15, 16).Value = 23750
$ws.Cells.Item(15, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(15, 19).Value = 1319
$ws.Cells.Item(15, 20).Value = 18

# Row 16
$ws.Cells.Item(16, 4).Value = 44181
$ws.Cells.Item(16, 11).Value = 'Modesto'
$ws.Cells.Item(16, 12).Value = 'Primera'
$ws.Cells.Item(16, 13).Value = 16
$ws.Cells.Item(16, 14).Value = 495000
$ws.Cells.Item(16, 15).Value = 500000
$ws.Cells.Item(16, 16).Value = 497500
$ws.Cells.Item(16, 17).Value = '$/bins (500 kilos)'
$ws.Cells.Item(16, 19).Value = 995
$ws.Cells.Item(16, 20).Value = 500

# Row 17
$ws.Cells.Item(17, 4).Value = 44181
$ws.Cells.Item(17, 11).Value = 'Modesto'
$ws.Cells.Item(17, 12).Value = 'Segunda'
$ws.Cells.Item(17, 13).Value = 10
$ws.Cells.Item(17, 14).Value = 425000
$ws.Cells.Item(17, 15).Value = 430000
$ws.Cells.Item(17, 16).Value = 427500
$ws.Cells.Item(17, 17).Value = '$/bins (500 kilos)'
$ws.Cells.Item(17, 19).Value = 855
$ws.Cells.Item(17, 20).Value = 500
